# #5: property building done
# Adds one "property" row to each of the five asset-type sheets
# (汽車/存款/股票/基金受益憑證/保險), matching the source scraper's output
# for 林岱樺's 2012-02-29 disclosure (new stock holding + re-emitted rows).

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Set-RowFormat($ws, [int]$row, [int]$lastCol) {
    # Column A on a data row carries the bold/bordered "index" style (same
    # style as every other existing data row's A cell); copy it from a
    # known-good source cell elsewhere in the workbook so the freshly
    # inserted row matches the look of its neighbours.
    $srcA = $wb.Worksheets.Item(2).Range("A3")
    $srcA.Copy()
    $ws.Cells.Item($row, 1).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# Sheet 1 - 汽車 (car): single existing row, append row 2 (copy of row 1)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = 29
$ws1.Range("B2").Value = $ws1.Range("B1").Value2
$ws1.Range("C2").Value = $ws1.Range("C1").Value2
$ws1.Range("D2").Value = $ws1.Range("D1").Value2
$ws1.Range("E2").Value = $ws1.Range("E1").Value2
$ws1.Range("F2").Value = $ws1.Range("F1").Value2
$ws1.Range("G2").Value = $ws1.Range("G1").Value2
Set-RowFormat $ws1 2 7

# ---------------------------------------------------------------------
# Sheet 2 - 存款 (deposit): insert a new row 2 (copy of old row 1),
# push the existing 6 rows down (their "index" column shifts -1).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Insert()
$ws2.Range("A2").Value = 43
$ws2.Range("B2").Value = $ws2.Range("B1").Value2
$ws2.Range("C2").Value = $ws2.Range("C1").Value2
$ws2.Range("D2").Value = $ws2.Range("D1").Value2
$ws2.Range("E2").Value = $ws2.Range("E1").Value2
$ws2.Range("G2").Value = $ws2.Range("G1").Value2
Set-RowFormat $ws2 2 7

$ws2.Range("A3").Value = 43
$ws2.Range("A4").Value = 44
$ws2.Range("A5").Value = 46
$ws2.Range("A6").Value = 47
$ws2.Range("A7").Value = 48
$ws2.Range("A8").Value = 49

# ---------------------------------------------------------------------
# Sheet 3 - 股票 (stock): single existing header-ish row, append row 2
# with the brand-new stock holding (new unique strings introduced here).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = 56
$ws3.Range("B2").Value = "中日國際企業股份有限公司"
$ws3.Range("C2").Value = $ws3.Range("K1").Value2
$ws3.Range("D2").Value = 500
$ws3.Range("E2").Value = 10
$ws3.Range("F2").Value = $ws3.Range("D1").Value2
$ws3.Range("G2").Value = 5000
$ws3.Range("H2").Value = "stock"
$ws3.Range("I2").Value = "normal"
$ws3.Range("J2").Value = "2012-02-29"
$ws3.Range("K2").Value = $ws3.Range("K1").Value2
$ws3.Range("L2").Value = 904
$ws3.Range("M2").Value = "tmp3bff1"
$ws3.Range("N2").Value = 56
Set-RowFormat $ws3 2 14

# ---------------------------------------------------------------------
# Sheet 4 - 基金受益憑證 (funds): insert a new row 2 (copy of old row 1),
# push the existing 2 rows down (their "index" column shifts -1).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Insert()
$ws4.Range("A2").Value = 66
$ws4.Range("B2").Value = $ws4.Range("B1").Value2
$ws4.Range("C2").Value = $ws4.Range("C1").Value2
$ws4.Range("D2").Value = $ws4.Range("D1").Value2
$ws4.Range("E2").Value = $ws4.Range("E1").Value2
$ws4.Range("F2").Value = $ws4.Range("F1").Value2
$ws4.Range("G2").Value = $ws4.Range("G1").Value2
$ws4.Range("H2").Value = $ws4.Range("H1").Value2
Set-RowFormat $ws4 2 8

$ws4.Range("A3").Value = 66
$ws4.Range("A4").Value = 67

# ---------------------------------------------------------------------
# Sheet 5 - 保險 (insurance): insert a new row 2 (copy of old row 1),
# push the existing 1 row down (its "index" column is unchanged, =81).
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows.Item(2).Insert()
$ws5.Range("A2").Value = 80
$ws5.Range("B2").Value = $ws5.Range("B1").Value2
$ws5.Range("C2").Value = $ws5.Range("C1").Value2
$ws5.Range("D2").Value = $ws5.Range("D1").Value2
Set-RowFormat $ws5 2 4

$ws5.Range("A3").Value = 81

Write-Output "done"
